# Commit: "Cartes -> Anneau de vérité"
# Adds a new "Triche" (cheat) sub-category plus several new cards
# (Maître de son destin, Espoir éternel, Capuche de parieur, Joyau du
# parieur, Anneau de Kallas, Anneau de concentration, Anneau de
# sacrifice, Anneau de vérité) to the card-repartition sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Touch the brand-new shared-string values first, in the exact order the
# author typed them, so the shared-strings table grows in the same
# sequence as the authored file (Nom column first, then the new
# "Triche" Catégorie 2 label, then the remaining new Nom values).
$ws.Cells.Item(3, 1).Value  = "Maître de son destin"    # A3
$ws.Cells.Item(55, 1).Value = "Espoir éternel"           # A55
$ws.Cells.Item(3, 10).Value = "Triche"                   # J3 (Catégorie 2)
$ws.Cells.Item(56, 1).Value = "Capuche de parieur"       # A56
$ws.Cells.Item(57, 1).Value = "Joyau du parieur"         # A57
$ws.Cells.Item(58, 1).Value = "Anneau de Kallas"         # A58
$ws.Cells.Item(59, 1).Value = "Anneau de concentration"  # A59
$ws.Cells.Item(60, 1).Value = "Anneau de sacrifice"      # A60
$ws.Cells.Item(61, 1).Value = "Anneau de vérité"         # A61

# --- Row 3 : Maître de son destin (Enchantement) ---
$ws.Cells.Item(3, 2).Value  = "Enchantement"             # B3
$ws.Cells.Item(3, 4).Value  = 1                          # D3 (Bleu)
$ws.Cells.Item(3, 6).Value  = 1                          # F3 (Rouge)
$ws.Cells.Item(3, 8).Value  = 1                          # H3 (Générique)

# --- Row 29 : register "Triche" label in the Catégorie-2 lookup table ---
$ws.Cells.Item(29, 20).Value = "Triche"                  # T29

# --- Row 53 : Espoir éternel becomes a "Triche" card instead of "Hasard" ---
$ws.Cells.Item(53, 10).Value = "Triche"                  # J53

# --- Row 55 : Espoir éternel (Enchantement) ---
$ws.Cells.Item(55, 2).Value  = "Enchantement"            # B55
$ws.Cells.Item(55, 3).Value  = 1                         # C55 (Blanc)
$ws.Cells.Item(55, 4).Value  = 1                         # D55 (Bleu)
$ws.Cells.Item(55, 8).Value  = 2                         # H55 (Générique)
$ws.Cells.Item(55, 10).Value = "Triche"                  # J55 (Catégorie 2)

# --- Row 56 : Capuche de parieur (Artefact) ---
$ws.Cells.Item(56, 2).Value  = "Artefact"                # B56
$ws.Cells.Item(56, 8).Value  = 2                         # H56 (Générique)
$ws.Cells.Item(56, 10).Value = "Triche"                  # J56
$ws.Cells.Item(56, 11).Value = "Hasard"                  # K56

# --- Row 57 : Joyau du parieur (Artefact) ---
$ws.Cells.Item(57, 2).Value  = "Artefact"                # B57
$ws.Cells.Item(57, 8).Value  = 2                         # H57
$ws.Cells.Item(57, 10).Value = "Triche"                  # J57
$ws.Cells.Item(57, 11).Value = "Ramp"                    # K57

# --- Row 58 : Anneau de Kallas (Artefact) ---
$ws.Cells.Item(58, 2).Value  = "Artefact"                # B58
$ws.Cells.Item(58, 8).Value  = 2                         # H58
$ws.Cells.Item(58, 10).Value = "Ramp"                    # J58

# --- Row 59 : Anneau de concentration (Artefact) ---
$ws.Cells.Item(59, 2).Value  = "Artefact"                # B59
$ws.Cells.Item(59, 8).Value  = 2                         # H59
$ws.Cells.Item(59, 10).Value = "Ramp"                    # J59

# --- Row 60 : Anneau de sacrifice (Artefact) ---
$ws.Cells.Item(60, 2).Value  = "Artefact"                # B60
$ws.Cells.Item(60, 8).Value  = 2                         # H60
$ws.Cells.Item(60, 10).Value = "Ramp"                    # J60

# --- Row 61 : Anneau de vérité (Artefact) ---
$ws.Cells.Item(61, 2).Value  = "Artefact"                # B61
$ws.Cells.Item(61, 8).Value  = 2                         # H61
$ws.Cells.Item(61, 10).Value = "Ramp"                    # J61

# --- Scroll / selection state, as left by the author ---
$ws.Application.ActiveWindow.ScrollRow = 43
$ws.Range("A62").Select() | Out-Null
